$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.376.09"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.614.88"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.487"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.840.33"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.614.61"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "26.378.10"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").Value = "1.165.11"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").Value = "1.753.59"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.406"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "0.0₇0956"
$ws.Range("E50").Value = "  -10.45%  "
$ws.Range("E51").Value = "  -0.05%  "
